$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price strings (e.g. "326.80", "6.200", "0.4420") must be
# written with a leading apostrophe so Excel keeps them as literal text
# instead of silently reparsing them as numbers and dropping trailing zeros.

$ws.Range("D2").Value = '27.723.87'
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").Value = '1.758.81'
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").Value = "'326.80"
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = "'0.4420"
$ws.Range("E7").Value = '  -2.21%  '
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").Value = "'45.85"
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D10").Value = "'0.07761"
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("D11").Value = "'1.128"
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").Value = "'21.78"
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("D14").Value = "'6.200"
$ws.Range("E14").Value = '  -1.77%  '
$ws.Range("D15").Value = "'7.364"
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").Value = '1.759.85'
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").Value = "'91.72"
$ws.Range("E17").Value = '  +13.03%  '
$ws.Range("D18").Value = "'0.00001081"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'0.06230"
$ws.Range("E19").Value = '  -7.56%  '
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = "'17.42"
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("D22").Value = "'6.200"
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("D23").Value = "'0.5342"
$ws.Range("E23").Value = '  -3.18%  '
$ws.Range("D24").Value = '27.764.77'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("D26").Value = "'2.335"
$ws.Range("E26").Value = '  -3.76%  '
$ws.Range("D27").Value = "'20.87"
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").Value = "'153.56"
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").Value = "'2.369"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '1.959.09'
$ws.Range("D31").Value = "'129.14"
$ws.Range("E31").Value = '  -3.14%  '
$ws.Range("D32").Value = "'1.216"
$ws.Range("E32").Value = '  -1.76%  '
$ws.Range("D33").Value = "'5.780"
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("D34").Value = "'0.09285"
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = "'3.669"
$ws.Range("E35").Value = '  -8.83%  '
$ws.Range("D36").Value = "'12.75"
$ws.Range("E36").Value = '  +4.67%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.02346"
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2194"
$ws.Range("E38").Value = '  -6.54%  '
$ws.Range("D39").Value = "'0.6542"
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("D40").Value = "'5.109"
$ws.Range("E40").Value = '  -1.86%  '
$ws.Range("D41").Value = "'0.06143"
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").Value = "'1.200"
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").Value = "'8.041"
$ws.Range("E43").Value = '  -3.99%  '
$ws.Range("D44").Value = "'1.417"
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = "'13.85"
$ws.Range("E46").Value = '  -1.94%  '
$ws.Range("D47").Value = "'0.6033"
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").Value = "'3.761"
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").Value = "'126.30"
$ws.Range("E49").Value = '  -2.93%  '
$ws.Range("D50").Value = "'2.001"
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").Value = "'1.148"
$ws.Range("E51").Value = '  -1.58%  '
